$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("day1").Range("B2").Value = 90
$wb.Worksheets.Item("day1").Range("B5").Value = 10

$wb.Worksheets.Item("day2").Range("B2").Value = 90
$wb.Worksheets.Item("day2").Range("B4").Value = 0
$wb.Worksheets.Item("day2").Range("B5").Value = 10

$wb.Worksheets.Item("day3").Range("B2").Value = 90
$wb.Worksheets.Item("day3").Range("B5").Value = 10

$wb.Worksheets.Item("day4").Range("B2").Value = 87.5
$wb.Worksheets.Item("day4").Range("B5").Value = 12.5
$wb.Worksheets.Item("day4").Range("B10").Value = 0

$wb.Worksheets.Item("day5").Range("B2").Value = 100
$wb.Worksheets.Item("day5").Range("B5").Value = 0
$wb.Worksheets.Item("day5").Range("B6").Value = 0
$wb.Worksheets.Item("day5").Range("B7").Value = 0
$wb.Worksheets.Item("day5").Range("B8").Value = 0
$wb.Worksheets.Item("day5").Range("B10").Value = 0

$wb.Worksheets.Item("day6").Range("B2").Value = 37.5
$wb.Worksheets.Item("day6").Range("B5").Value = 50
$wb.Worksheets.Item("day6").Range("B10").Value = 12.5
